# Generate Report for Handoff
# Update the report UUID/token and timestamps across all sheets to reflect
# a newly generated handoff report.

$wb = $excel.ActiveWorkbook

$oldGuid = "935889ee-e818-4893-8ccb-f3a6808c91ad"
$newGuid = "78d1c34f-5577-43e5-b418-0e1ad8a791da"

$oldZhHash = "0e7d5c93e6885fda40509d1937d9ccccd3d3e06e"
$newZhHash = "eae86f19b4e49edd48b23ab17e1e8c13bd63c19a"

$newMdName = "$newGuid.md"
$newZhXlf = "$newGuid.$newZhHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newZhHash.de-de.xlf"

$newHandoffDate = "2016-03-23 05:08:01"
$newZhHandoffDate = "2016-03-23 05:07:57"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newHandoffDate
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newMdName

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhHandoffDate
$wsZh.Hyperlinks.Item(1).TextToDisplay = $newMdName
$wsZh.Hyperlinks.Item(2).TextToDisplay = $newZhXlf

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Hyperlinks.Item(1).TextToDisplay = $newMdName
$wsDe.Hyperlinks.Item(2).TextToDisplay = $newDeXlf
